# Updates cryptos list values (price + 1h volume change) per the
# Wed Oct 23 15:28:26 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.219.17"
$ws.Range("E2").Value = "  -1.99%  "

$ws.Range("D3").Value = "2.544.82"
$ws.Range("E3").Value = "  -3.51%  "

$ws.Range("E4").Value = "  +0.72%  "

$ws.Range("D5").Value = "'579.94"
$ws.Range("E5").Value = "  -2.89%  "

$ws.Range("D6").Value = "'166.85"
$ws.Range("E6").Value = "  -2.06%  "

$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  -1.83%  "

$ws.Range("D9").Value = "2.544.30"
$ws.Range("E9").Value = "  -3.31%  "

$ws.Range("E10").Value = "  -1.45%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "'0.351"
$ws.Range("E12").Value = "  -3.68%  "

$ws.Range("D13").Value = "'5.12"
$ws.Range("E13").Value = "  -2.63%  "

$ws.Range("D14").Value = "'26.46"
$ws.Range("E14").Value = "  -4.89%  "

$ws.Range("D15").Value = "3.013.15"
$ws.Range("E15").Value = "  -2.82%  "

$ws.Range("D16").Value = "'0.0000176"
$ws.Range("E16").Value = "  -3.52%  "

$ws.Range("D17").Value = "66.136.21"
$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("D18").Value = "2.544.35"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").Value = "'11.31"
$ws.Range("E19").Value = "  -6.98%  "

$ws.Range("D20").Value = "'7.65"
$ws.Range("E20").Value = "  -5.98%  "

$ws.Range("D21").Value = "'347.09"
$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").Value = "'4.19"
$ws.Range("E22").Value = "  -3.58%  "

$ws.Range("D23").Value = "'4.54"
$ws.Range("E23").Value = "  -3.48%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").Value = "'1.90"
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").Value = "'68.94"
$ws.Range("E26").Value = "  -1.57%  "

$ws.Range("D27").Value = "'9.84"
$ws.Range("E27").Value = "  -7.40%  "

$ws.Range("D28").Value = "2.684.37"
$ws.Range("E28").Value = "  -2.70%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").Value = "0.0₃0974"
$ws.Range("E30").Value = "  -3.65%  "

$ws.Range("D31").Value = "'8.18"
$ws.Range("E31").Value = "  +2.71%  "

$ws.Range("D32").Value = "'522.37"
$ws.Range("E32").Value = "  -5.70%  "

$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "  -3.97%  "

$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -5.22%  "

$ws.Range("D35").Value = "'0.130"
$ws.Range("E35").Value = "  -5.05%  "

$ws.Range("E36").Value = "  +0.81%  "

$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").Value = "'1.44"
$ws.Range("E38").Value = "  -4.38%  "

$ws.Range("D39").Value = "'18.64"
$ws.Range("E39").Value = "  -2.17%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'18.28"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.354"
$ws.Range("E41").Value = "  -3.82%  "

$ws.Range("D42").Value = "'1.76"
$ws.Range("E42").Value = "  -3.15%  "

$ws.Range("D43").Value = "'5.05"
$ws.Range("E43").Value = "  -3.40%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'2.39"
$ws.Range("E45").Value = "  -1.66%  "

$ws.Range("D46").Value = "0.0₆0282"
$ws.Range("E46").Value = "  -5.40%  "

$ws.Range("D47").Value = "'147.58"
$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("D48").Value = "'0.557"
$ws.Range("E48").Value = "  -4.52%  "

$ws.Range("D49").Value = "'3.68"
$ws.Range("E49").Value = "  -3.14%  "

$ws.Range("D50").Value = "'1.70"
$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("D51").Value = "'0.0758"
$ws.Range("E51").Value = "  -1.94%  "
